# proyecto.xlsx - correccion de errores en datos (ventas/inventario/clientes)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja "Clientes": se eliminan los dos ultimos registros (filas 6 y 7)
# ---------------------------------------------------------------------------
$wsClientes = $wb.Worksheets.Item("Clientes")
$wsClientes.Rows("6:7").Delete()

# ---------------------------------------------------------------------------
# Hoja "Inventario": se corrigen existencias de varios productos
# ---------------------------------------------------------------------------
$wsInventario = $wb.Worksheets.Item("Inventario")

$wsInventario.Range("C2").NumberFormat = "@"
$wsInventario.Range("C2").Value = "315"   # Coca Cola: 310 -> 315
$wsInventario.Range("E2").NumberFormat = "@"
$wsInventario.Range("E2").Value = "5"

$wsInventario.Range("C3").Value = 245     # Pepsi 500 ml: 230 -> 245

$wsInventario.Range("C4").NumberFormat = "@"
$wsInventario.Range("C4").Value = "400"   # Monster Ultra White: 391 -> 400
$wsInventario.Range("E4").NumberFormat = "@"
$wsInventario.Range("E4").Value = "18.5"

$wsInventario.Range("C5").NumberFormat = "@"
$wsInventario.Range("C5").Value = "450"
$wsInventario.Range("E5").NumberFormat = "@"
$wsInventario.Range("E5").Value = "19"

$wsInventario.Range("C6").Value = 87      # Aritos: 92 -> 87

# ---------------------------------------------------------------------------
# Hoja "Ventas": se corrige el primer registro y se eliminan el resto de
# transacciones ya procesadas
# ---------------------------------------------------------------------------
$wsVentas = $wb.Worksheets.Item("Ventas")
$wsVentas.Range("A2").NumberFormat = "@"
$wsVentas.Range("A2").Value = "251"
$wsVentas.Range("C2").Value = 5
$wsVentas.Range("D2").Value = 10
$wsVentas.Rows("3:11").Delete()

[void]$wsVentas.Range("C14").Select()
